# Generate Report for Handoff
# Regenerates handoff data for the zh-cn and de-de locale sheets:
#   - Priority for every "Ready for handoff" row moves from "ht" to "mt"
#   - Latest Handoff Datetime is refreshed to a newer timestamp
# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff datetime, so it is refreshed to match as well.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

for ($row = 3; $row -le 16; $row++) {
    $wsZhCn.Range("E$row").Value = "mt"
    $wsZhCn.Range("H$row").Value = "2016-08-22 02:34:32"

    $wsDeDe.Range("E$row").Value = "mt"
    $wsDeDe.Range("H$row").Value = "2016-08-22 02:34:38"

    $wsOverview.Range("G$row").Value = "2016-08-22 02:34:38"
}
